# Update countries & provincias Spain
# - Swap order of "Japon"/"Polonia" in the country list (row 47 becomes
#   Polonia with updated stats, row 48 becomes Japon keeping the old
#   Polonia-position stats that previously belonged to Japon).
# - Refresh the "Datos actualizados..." timestamp cell.
# - Update the numeric stats for several countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swap (row 47 / row 48) -----------------------------
$ws.Range("A47").Value = "Polonia"
$ws.Range("A48").Value = "Japon"

# --- Timestamp update --------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 11:01"

# --- Numeric data updates ----------------------------------------------
# Row => column letter => new value
$updates = @{
    7   = @{ B = 995319; C = 4993; D = 809387; E = 168756; G = 83; H = 17176 }
    25  = @{ B = 220819; C = 3446; D = 157562; E = 59699;  G = 38; H = 3558 }
    26  = @{ B = 174796; C = 2743; D = 125959; E = 41420;  G = 74; H = 7417 }
    31  = @{ B = 115057; C = 1037; D = 93801;  E = 20334;  G = 3;  H = 922 }
    47  = @{ B = 67372;  C = 502;  D = 46638;  E = 18695;  G = 6;  H = 2039 }
    48  = @{ B = 67264;  D = 56164; E = 9836;  H = 1264 }
    52  = @{ B = 56812;  C = 41;   E = 1199 }
    65  = @{ B = 38165;  C = 3;    E = 7674 }
    96  = @{ B = 9340;   C = 6;    D = 9054;  E = 159;   G = 1; H = 127 }
    98  = @{ B = 8550;   D = 7341; E = 1141 }
    111 = @{ B = 4811;   C = 9;    D = 4342;  E = 380 }
    120 = @{ B = 3917;   C = 41;   D = 2377;  E = 1507 }
    130 = @{ B = 2906;   C = 32;   D = 1840;  E = 980 }
    131 = @{ B = 2883;   C = 18;   D = 2296;  E = 454 }
    136 = @{ B = 2375;   C = 2;    E = 223 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
